# Plantilla de usuarios - carga masiva
# - renombra la hoja
# - quita la fila de ejemplo (D2) que llevaba el formato de fecha
# - da estilo de encabezado (negrita, bordes, centrado) a la fila 1
# - ajusta margenes de pagina y la celda seleccionada

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renombrar la hoja
$ws.Name = "Sheet1"

# Quitar la fila de ejemplo (fila 2), que solo contenia el formato de fecha sin datos
$ws.Range("A2:I2").EntireRow.Delete()

# Dar formato de encabezado a la fila de titulos: negrita, borde fino y centrado
$header = $ws.Range("A1:I1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# Margenes de pagina por defecto
$ws.PageSetup.LeftMargin = 0.7 * 72
$ws.PageSetup.RightMargin = 0.7 * 72
$ws.PageSetup.TopMargin = 0.75 * 72
$ws.PageSetup.BottomMargin = 0.75 * 72
$ws.PageSetup.HeaderMargin = 0.3 * 72
$ws.PageSetup.FooterMargin = 0.3 * 72

# Seleccion activa
$ws.Range("L7").Select() | Out-Null
